# Apply targeted odds updates to sheet1 (active sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 3.6  # G2: 3.8 -> 3.6
$ws.Cells.Item(2, 9).Value = 2.1  # I2: 2 -> 2.1
$ws.Cells.Item(2, 18).Value = 1.8  # R2: 1.91 -> 1.8
$ws.Cells.Item(2, 19).Value = 1.95  # S2: 1.91 -> 1.95
$ws.Cells.Item(2, 21).Value = 17  # U2: 19 -> 17
$ws.Cells.Item(2, 32).Value = 9.5  # AF2: 9 -> 9.5
$ws.Cells.Item(2, 34).Value = 19  # AH2: 17 -> 19
# Row 4
$ws.Cells.Item(4, 8).Value = 5.5  # H4: 4.75 -> 5.5
$ws.Cells.Item(4, 9).Value = 1.36  # I4: 1.3 -> 1.36
$ws.Cells.Item(4, 11).Value = 21  # K4: 19 -> 21
$ws.Cells.Item(4, 18).Value = 1.67  # R4: 1.7 -> 1.67
$ws.Cells.Item(4, 19).Value = 2.1  # S4: 2.05 -> 2.1
$ws.Cells.Item(4, 20).Value = 21  # T4: 26 -> 21
$ws.Cells.Item(4, 22).Value = 19  # V4: 23 -> 19
$ws.Cells.Item(4, 23).Value = 67  # W4: 81 -> 67
$ws.Cells.Item(4, 24).Value = 41  # X4: 51 -> 41
$ws.Cells.Item(4, 26).Value = 21  # Z4: 19 -> 21
$ws.Cells.Item(4, 27).Value = 11  # AA4: 10 -> 11
$ws.Cells.Item(4, 29).Value = 41  # AC4: 51 -> 41
$ws.Cells.Item(4, 31).Value = 10  # AE4: 9.5 -> 10
$ws.Cells.Item(4, 32).Value = 8.5  # AF4: 8 -> 8.5
$ws.Cells.Item(4, 34).Value = 10  # AH4: 9.5 -> 10
$ws.Cells.Item(4, 35).Value = 11  # AI4: 10 -> 11
# Row 5
$ws.Cells.Item(5, 10).Value = 1.17  # J5: 1.14 -> 1.17
$ws.Cells.Item(5, 11).Value = 5  # K5: 5.5 -> 5
# Row 6
$ws.Cells.Item(6, 7).Value = 3.6  # G6: 3.5 -> 3.6
$ws.Cells.Item(6, 9).Value = 2.25  # I6: 2.3 -> 2.25
$ws.Cells.Item(6, 10).Value = 1.11  # J6: 1.1 -> 1.11
$ws.Cells.Item(6, 11).Value = 6.5  # K6: 7 -> 6.5
$ws.Cells.Item(6, 21).Value = 17  # U6: 15 -> 17
$ws.Cells.Item(6, 27).Value = 6  # AA6: 5.5 -> 6
$ws.Cells.Item(6, 32).Value = 9.5  # AF6: 10 -> 9.5
$ws.Cells.Item(6, 35).Value = 21  # AI6: 23 -> 21
# Row 10
$ws.Cells.Item(10, 7).Value = 2.3  # G10: 2.25 -> 2.3
$ws.Cells.Item(10, 9).Value = 2.8  # I10: 2.88 -> 2.8
$ws.Cells.Item(10, 12).Value = 1.18  # L10: 1.2 -> 1.18
$ws.Cells.Item(10, 13).Value = 4.5  # M10: 4.33 -> 4.5
$ws.Cells.Item(10, 14).Value = 1.65  # N10: 1.67 -> 1.65
$ws.Cells.Item(10, 15).Value = 2.2  # O10: 2.15 -> 2.2
$ws.Cells.Item(10, 25).Value = 23  # Y10: 21 -> 23
$ws.Cells.Item(10, 29).Value = 41  # AC10: 34 -> 41
$ws.Cells.Item(10, 36).Value = 23  # AJ10: 26 -> 23
# Row 12
$ws.Cells.Item(12, 7).Value = 1.8  # G12: 1.72 -> 1.8
$ws.Cells.Item(12, 9).Value = 4.6  # I12: 5 -> 4.6
$ws.Cells.Item(12, 18).Value = 2.18  # R12: 2.22 -> 2.18
$ws.Cells.Item(12, 19).Value = 1.53  # S12: 1.52 -> 1.53
$ws.Cells.Item(12, 20).Value = 5  # T12: 4.9 -> 5
$ws.Cells.Item(12, 21).Value = 7  # U12: 6.7 -> 7
$ws.Cells.Item(12, 23).Value = 14  # W12: 13 -> 14
$ws.Cells.Item(12, 24).Value = 18.5  # X12: 17.5 -> 18.5
$ws.Cells.Item(12, 26).Value = 6.6  # Z12: 6.5 -> 6.6
$ws.Cells.Item(12, 31).Value = 9.5  # AE12: 10 -> 9.5
$ws.Cells.Item(12, 32).Value = 24  # AF12: 27 -> 24
$ws.Cells.Item(12, 33).Value = 16.5  # AG12: 18 -> 16.5
$ws.Cells.Item(12, 34).Value = 80  # AH12: 100 -> 80
$ws.Cells.Item(12, 35).Value = 60  # AI12: 70 -> 60
$ws.Cells.Item(12, 36).Value = 80  # AJ12: 90 -> 80
# Row 13
$ws.Cells.Item(13, 7).Value = 2.15  # G13: 2.2 -> 2.15
$ws.Cells.Item(13, 8).Value = 3.3  # H13: 3.25 -> 3.3
$ws.Cells.Item(13, 9).Value = 3.4  # I13: 3.3 -> 3.4
$ws.Cells.Item(13, 10).Value = 1.08  # J13: 1.07 -> 1.08
$ws.Cells.Item(13, 11).Value = 8  # K13: 9 -> 8
$ws.Cells.Item(13, 18).Value = 2  # R13: 1.95 -> 2
$ws.Cells.Item(13, 19).Value = 1.75  # S13: 1.8 -> 1.75
$ws.Cells.Item(13, 23).Value = 19  # W13: 21 -> 19
$ws.Cells.Item(13, 32).Value = 17  # AF13: 15 -> 17
$ws.Cells.Item(13, 33).Value = 13  # AG13: 12 -> 13
$ws.Cells.Item(13, 35).Value = 34  # AI13: 29 -> 34
# Row 17
$ws.Cells.Item(17, 7).Value = 2.15  # G17: 2.1 -> 2.15
$ws.Cells.Item(17, 9).Value = 3.7  # I17: 3.9 -> 3.7
$ws.Cells.Item(17, 21).Value = 9.5  # U17: 9 -> 9.5
$ws.Cells.Item(17, 22).Value = 9.5  # V17: 9 -> 9.5
$ws.Cells.Item(17, 32).Value = 17  # AF17: 19 -> 17
# Row 19
$ws.Cells.Item(19, 7).Value = 1.95  # G19: 2 -> 1.95
$ws.Cells.Item(19, 8).Value = 3.7  # H19: 3.6 -> 3.7
$ws.Cells.Item(19, 9).Value = 3.6  # I19: 3.4 -> 3.6
$ws.Cells.Item(19, 22).Value = 8.5  # V19: 9 -> 8.5
$ws.Cells.Item(19, 33).Value = 13  # AG19: 12 -> 13
$ws.Cells.Item(19, 35).Value = 29  # AI19: 26 -> 29
$ws.Cells.Item(19, 36).Value = 34  # AJ19: 29 -> 34
# Row 20
$ws.Cells.Item(20, 21).Value = 15.5  # U20: 16 -> 15.5
$ws.Cells.Item(20, 24).Value = 35  # X20: 32 -> 35
$ws.Cells.Item(20, 25).Value = 45  # Y20: 50 -> 45
$ws.Cells.Item(20, 27).Value = 6.3  # AA20: 6.2 -> 6.3
$ws.Cells.Item(20, 28).Value = 17.5  # AB20: 17 -> 17.5
$ws.Cells.Item(20, 31).Value = 6.3  # AE20: 6.2 -> 6.3
# Row 21
$ws.Cells.Item(21, 9).Value = 3.05  # I21: 3.1 -> 3.05
$ws.Cells.Item(21, 18).Value = 1.87  # R21: 1.85 -> 1.87
$ws.Cells.Item(21, 26).Value = 7.4  # Z21: 7.3 -> 7.4
$ws.Cells.Item(21, 27).Value = 5.9  # AA21: 5.8 -> 5.9
# Row 23
$ws.Cells.Item(23, 14).Value = 1.94  # N23: 2 -> 1.94
$ws.Cells.Item(23, 15).Value = 1.74  # O23: 1.8 -> 1.74
# Row 24
$ws.Cells.Item(24, 7).Value = 3.4  # G24: 3.5 -> 3.4
$ws.Cells.Item(24, 14).Value = 1.79  # N24: 1.85 -> 1.79
$ws.Cells.Item(24, 15).Value = 1.89  # O24: 1.95 -> 1.89
# Row 25
$ws.Cells.Item(25, 7).Value = 1.25  # G25: 1.3 -> 1.25
$ws.Cells.Item(25, 8).Value = 6  # H25: 5.75 -> 6
$ws.Cells.Item(25, 9).Value = 8  # I25: 7 -> 8
$ws.Cells.Item(25, 23).Value = 8  # W25: 8.5 -> 8
$ws.Cells.Item(25, 26).Value = 19  # Z25: 17 -> 19
$ws.Cells.Item(25, 28).Value = 23  # AB25: 21 -> 23
$ws.Cells.Item(25, 29).Value = 67  # AC25: 51 -> 67
$ws.Cells.Item(25, 31).Value = 23  # AE25: 21 -> 23
$ws.Cells.Item(25, 33).Value = 26  # AG25: 21 -> 26
$ws.Cells.Item(25, 34).Value = 101  # AH25: 81 -> 101
$ws.Cells.Item(25, 36).Value = 51  # AJ25: 41 -> 51
# Row 26
$ws.Cells.Item(26, 7).Value = 3.7  # G26: 3.6 -> 3.7
$ws.Cells.Item(26, 9).Value = 1.75  # I26: 1.8 -> 1.75
$ws.Cells.Item(26, 27).Value = 8.5  # AA26: 8 -> 8.5
$ws.Cells.Item(26, 31).Value = 11  # AE26: 10 -> 11
# Row 27
$ws.Cells.Item(27, 7).Value = 9  # G27: 8.5 -> 9
$ws.Cells.Item(27, 9).Value = 1.2  # I27: 1.22 -> 1.2
$ws.Cells.Item(27, 10).Value = 26  # J27: 1.01 -> 26
$ws.Cells.Item(27, 11).Value = 1.02  # K27: 15 -> 1.02
$ws.Cells.Item(27, 18).Value = 1.8  # R27: 1.83 -> 1.8
$ws.Cells.Item(27, 19).Value = 1.91  # S27: 1.83 -> 1.91
$ws.Cells.Item(27, 25).Value = 51  # Y27: 41 -> 51
$ws.Cells.Item(27, 26).Value = 26  # Z27: 23 -> 26
# Row 29
$ws.Cells.Item(29, 7).Value = 2.4  # G29: 2.35 -> 2.4
$ws.Cells.Item(29, 9).Value = 2.8  # I29: 2.9 -> 2.8
$ws.Cells.Item(29, 10).Value = 1.07  # J29: 1.06 -> 1.07
$ws.Cells.Item(29, 11).Value = 9  # K29: 10 -> 9
$ws.Cells.Item(29, 18).Value = 1.91  # R29: 1.95 -> 1.91
$ws.Cells.Item(29, 19).Value = 1.91  # S29: 1.8 -> 1.91
$ws.Cells.Item(29, 20).Value = 7.5  # T29: 7 -> 7.5
$ws.Cells.Item(29, 26).Value = 9  # Z29: 8.5 -> 9
$ws.Cells.Item(29, 28).Value = 15  # AB29: 17 -> 15
$ws.Cells.Item(29, 30).Value = 301  # AD29: 351 -> 301
$ws.Cells.Item(29, 34).Value = 29  # AH29: 34 -> 29
$ws.Cells.Item(29, 35).Value = 23  # AI29: 26 -> 23
# Row 30
$ws.Cells.Item(30, 9).Value = 4  # I30: 4.1 -> 4
$ws.Cells.Item(30, 20).Value = 9.5  # T30: 9 -> 9.5
$ws.Cells.Item(30, 27).Value = 7.5  # AA30: 8 -> 7.5
# Row 32
$ws.Cells.Item(32, 7).Value = 1.8  # G32: 1.85 -> 1.8
$ws.Cells.Item(32, 8).Value = 3.75  # H32: 3.7 -> 3.75
$ws.Cells.Item(32, 9).Value = 4  # I32: 3.9 -> 4
$ws.Cells.Item(32, 24).Value = 13  # X32: 15 -> 13
$ws.Cells.Item(32, 27).Value = 7.5  # AA32: 7 -> 7.5
# Row 34
$ws.Cells.Item(34, 14).Value = 1.65  # N34: 1.62 -> 1.65
$ws.Cells.Item(34, 15).Value = 2.2  # O34: 2.25 -> 2.2
# Row 37
$ws.Cells.Item(37, 11).Value = 13  # K37: 12 -> 13
